$wb = $excel.ActiveWorkbook

# --- 1. Add new enum values MRM and PRM to the "data_collection_mode list" sheet ---
$listSheet = $wb.Worksheets.Item("data_collection_mode list")
$listSheet.Range("A3").Value = "MRM"
$listSheet.Range("A4").Value = "PRM"

# --- 2. Update the comment on the data_collection_mode column header (Y1) on the main sheet ---
$mainSheet = $wb.Worksheets.Item("Export as TSV")
$newComment = "Mode of data collection in tandem MS assays. Either DDA (Data-dependent acquisition), DIA (Data-independent acquisition), MRM (multiple reaction monitoring), or PRM (parallel reaction monitoring)."
$mainSheet.Range("Y1").Comment.Text($newComment)

# --- 3. Update the data validation on column Y to reference the expanded list and new error text ---
$dvRange = $mainSheet.Range("Y2:Y1048576")
$validation = $dvRange.Validation
$validation.Modify(
    [Microsoft.Office.Interop.Excel.XlDVType]::xlValidateList,
    [Microsoft.Office.Interop.Excel.XlFormatConditionOperator]::xlBetween,
    [Microsoft.Office.Interop.Excel.XlDVAlertStyle]::xlValidAlertStop,
    "'data_collection_mode list'!`$A`$1:`$A`$4"
)
$validation.ErrorTitle = "Value must come from list"
$validation.ErrorMessage = "Value must be one of: DDA / DIA / MRM / PRM."
